# Weekly update: insert the latest daily price record at the top of the
# data table (row 2), pushing all existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the first data row.
$ws.Rows(2).Insert()

# The inserted row inherits the header row's formatting; strip that so it
# matches the look of an ordinary data row, then re-apply the date format
# used by the other rows' "Fecha" column (D).
$ws.Range("A2:T2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new record.
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 45092
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107001
$ws.Range("J2").Value = "Caqui"
$ws.Range("K2").Value = "Mankaki"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 140
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18429
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 1024
$ws.Range("T2").Value = 18
